$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values (price / 1h volume change) per row,
# matching the refreshed coinranking.com scrape.
$updates = @{
    "D2" = "321.33"
    "E2" = "-3.36%"
    "D3" = "42.85"
    "E3" = "-6.05%"
    "D4" = "5.227"
    "E4" = "-4.94%"
    "D5" = "0.08238"
    "E5" = "-3.50%"
    "D6" = "4.320"
    "E6" = "-2.88%"
    "D7" = "1.790"
    "E7" = "-12.94%"
    "D8" = "0.9490"
    "E8" = "-3.90%"
    "D9" = "0.1122"
    "E9" = "-2.98%"
    "D10" = "0.1884"
    "E10" = "-1.76%"
    "D11" = "0.09413"
    "E11" = "-4.04%"
    "D12" = "0.04622"
    "E12" = "-2.07%"
    "D13" = "7.446"
    "E13" = "-21.57%"
    "E14" = "-0.18%"
    "D15" = "0.001307"
    "D16" = "0.005786"
    "E16" = "-2.99%"
    "E17" = "-0.67%"
    "E18" = "0.27%"
    "D19" = "0.3367"
    "E19" = "0.35%"
    "D20" = "0.1387"
    "E20" = "0.97%"
    "E21" = "-0.16%"
    "D22" = "0.04163"
    "E22" = "0.49%"
    "D23" = "0.001248"
    "E23" = "-4.22%"
    "D24" = "0.004283"
    "E24" = "-6.21%"
    "E25" = "-6.37%"
    "E26" = "-0.38%"
    "E38" = "-3.24%"
    "D39" = "0.05601"
    "E39" = "-2.36%"
    "D40" = "0.008168"
    "E40" = "3.31%"
    "D41" = "0.1404"
    "E41" = "-2.06%"
    "D42" = "0.006541"
    "E42" = "-9.97%"
    "D43" = "0.002116"
    "E43" = "0.17%"
    "D44" = "0.007666"
    "E44" = "-13.17%"
    "D45" = "0.3483"
    "E45" = "-1.87%"
    "D46" = "0.00006742"
    "E46" = "-4.86%"
    "E47" = "-0.38%"
    "D48" = "0.003069"
    "E48" = "-11.26%"
    "D49" = "0.004096"
    "E49" = "15.69%"
    "E50" = "-0.38%"
    "E51" = "-0.38%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text storage (the source column holds formatted
    # price/percent strings, not numeric values).
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
